$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.641.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.794.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.557"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.297"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0695"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.813.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.636"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.572.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0802"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "165.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("E28").Value = "  +2.46%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.87%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0522"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.422.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("E36").Value = "  +6.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.676"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.31%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0192"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "85.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.28%  "
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.933"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("E45").Value = "  +3.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("E51").Value = "  -5.04%  "
